$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "246.30"
$r.Style = "Normal"

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "5.474"
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.05636"
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "6.462"
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.8060"
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "1.044"
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.1439"
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.07400"
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.03187"
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.09269"
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.001666"
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "3.204"
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.04736"
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.0005838"
$r.Style = "Normal"

$ws.Range("E17").Value = "16OneONE"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.006294"
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.001060"
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.004117"
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "3.981"
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.142"
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.1315"
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.0003009"
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.006881"
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.1039"
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.009041"
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.00005663"
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.6821"
$r.Style = "Normal"

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.01989"
$r.Style = "Normal"

$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.00002106"
$r.Style = "Normal"
